# Regenerate the lattice-multiplication exercise table: every cell keeps its
# original 5-line layout (problem, digit header, dashes, two lattice rows)
# but the numbers themselves are replaced with a newly generated set.
# Cell [r,c] content is replaced in place (table stays 5 rows x 3 cols).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Word represents a manual line break (<w:br/>) as chr(11) inside Range.Text.
$nl = [char]11

$newCells = @(
    @(1,1,"77 x 76","  7    6","  ----","7|    |","7|    |"),
    @(1,2,"67 x 91","  9    1","  ----","6|    |","7|    |"),
    @(1,3,"61 x 51","  5    1","  ----","6|    |","1|    |"),

    @(2,1,"48 x 22","  2    2","  ----","4|    |","8|    |"),
    @(2,2,"86 x 46","  4    6","  ----","8|    |","6|    |"),
    @(2,3,"84 x 79","  7    9","  ----","8|    |","4|    |"),

    @(3,1,"17 x 71","  7    1","  ----","1|    |","7|    |"),
    @(3,2,"21 x 83","  8    3","  ----","2|    |","1|    |"),
    @(3,3,"31 x 31","  3    1","  ----","3|    |","1|    |"),

    @(4,1,"23 x 39","  3    9","  ----","2|    |","3|    |"),
    @(4,2,"59 x 60","  6    0","  ----","5|    |","9|    |"),
    @(4,3,"66 x 46","  4    6","  ----","6|    |","6|    |"),

    @(5,1,"82 x 14","  1    4","  ----","8|    |","2|    |"),
    @(5,2,"33 x 22","  2    2","  ----","3|    |","3|    |"),
    @(5,3,"14 x 46","  4    6","  ----","1|    |","4|    |")
)

foreach ($entry in $newCells) {
    $row = $entry[0]
    $col = $entry[1]
    $lines = $entry[2..6]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = [string]::Join($nl, $lines)
}
